$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "1.000 GHz"
$ws.Range("C1").Value = "1234.000 K"

$ws.Range("A2").Value = "1.500 GHz"

$ws.Range("C3").Value = "1217.784 K"
$ws.Range("D3").Value = "1.427dB"

$ws.Range("B4").Value = "7.173 dB"

$ws.Range("B5").Value = "7.245 d8"
$ws.Range("C5").Value = "1247.628 K"

$ws.Range("A6").Value = "3.500GHz"
$ws.Range("B6").Value = "7.451 d8"

$ws.Range("A7").Value = "4.000 GHz"
$ws.Range("C7").Value = "1381.940 K"
